$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the previously-empty GWPb (E) and GWP-LULUC (F) columns with 0
# for the transportation rows (2-5), matching the updated background data.
$ws.Range("E2:F5").Value = 0

# Update the active selection to reflect where the user last clicked.
$ws.Range("C8").Select()
